$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Tahun), shifting NIK/Nama Karyawan/Alpa/Ijin/Terlambat right by one.
$ws.Columns("B").Insert()

# Update header row text
$ws.Range("A1").Value = "Bulan"
$ws.Range("B1").Value = "Tahun"

# A1:C1 (Bulan, Tahun, NIK) should be highlighted red
$ws.Range("A1:C1").Interior.Color = 255

# D1:G1 (Nama Karyawan, Alpa, Ijin, Terlambat) should be light blue
$ws.Range("D1:G1").Interior.Color = 15722206

# Fill in first data row (row 2)
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 2024
$ws.Range("C2").Value = "EN-4-034"
$ws.Range("D2").Value = "Nurdiansah"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Update selection to A3 to match author's final cursor position
$ws.Range("A3").Select() | Out-Null
